# "Generate Report for Handback"
#
# Two files (38caa61e-3a66-441d-97cc-1d7377981442 and
# 76a5350b-de9e-4e6b-acfd-af5a4cdc3ec6) that were previously "Ready for
# handoff" have now been handed back (translated) for both target
# languages (zh-cn and de-de). This script updates the localization
# status report to reflect that handback:
#   - Overview sheet: status columns for zh-cn/de-de flip to
#     "Handed back: in sync with en-US" for those two files.
#   - zh-cn / de-de detail sheets: Status flips the same way, the
#     "Latest Target File" / "Latest Handback File" columns are filled
#     in (previously blank) and a "Latest Handback DateTime" timestamp
#     is recorded (replacing the "0001-01-01 00:00:00" placeholder).
#   - New hyperlinks are added on the newly-populated "Latest Target
#     File" cells, mirroring the existing ones for the rows that were
#     already handed back.

$wb = $excel.ActiveWorkbook

$handedBack = "Handed back: in sync with en-US"

# ---------------------------------------------------------------------
# Overview sheet: rows 4 and 5 (the two files moving from "Ready for
# handoff" to handed-back) get their zh-cn (E) and de-de (F) status
# columns updated. Column G (Latest HO Xliff Generate Date) is
# unchanged.
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("E4").Value = $handedBack
$wsOverview.Range("F4").Value = $handedBack
$wsOverview.Range("E5").Value = $handedBack
$wsOverview.Range("F5").Value = $handedBack

# ---------------------------------------------------------------------
# zh-cn sheet: rows 4 and 5.
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$zhCnRows = @(
    @{ Row = 4; File = "38caa61e-3a66-441d-97cc-1d7377981442.md"; TargetRepo = "ol-test0-zhcn"; TargetCommit = "e10bb4d4424a3e8bed2f53b9c0170729a6741d33"; Path = "e2e/38caa61e-3a66-441d-97cc-1d7377981442.md" },
    @{ Row = 5; File = "76a5350b-de9e-4e6b-acfd-af5a4cdc3ec6.md"; TargetRepo = "ol-test0-zhcn"; TargetCommit = "e10bb4d4424a3e8bed2f53b9c0170729a6741d33"; Path = "e2e/76a5350b-de9e-4e6b-acfd-af5a4cdc3ec6.md" }
)

foreach ($info in $zhCnRows) {
    $r = $info.Row
    $wsZhCn.Cells.Item($r, 3).Value = $handedBack               # C: Status
    $xlf = $wsZhCn.Cells.Item($r, 7).Value()                    # G: Latest Handoff File
    $wsZhCn.Cells.Item($r, 9).Value = $info.File                # I: Latest Target File
    $wsZhCn.Cells.Item($r, 10).Value = $xlf                     # J: Latest Handback File
    $wsZhCn.Cells.Item($r, 11).Value = "2016-08-17 06:26:27"    # K: Latest Handback DateTime

    $target = "https://github.com/OpenLocalizationTestOrg/" + $info.TargetRepo + "/blob/" + $info.TargetCommit + "/" + $info.Path
    $wsZhCn.Hyperlinks.Add($wsZhCn.Cells.Item($r, 9), $target, "", "", $info.File)
}

# ---------------------------------------------------------------------
# de-de sheet: rows 4 and 5.
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$deDeRows = @(
    @{ Row = 4; File = "38caa61e-3a66-441d-97cc-1d7377981442.md"; TargetRepo = "ol-test0-dede"; TargetCommit = "be90939c9de779fd82751107be1eb50673a8d746"; Path = "e2e/38caa61e-3a66-441d-97cc-1d7377981442.md" },
    @{ Row = 5; File = "76a5350b-de9e-4e6b-acfd-af5a4cdc3ec6.md"; TargetRepo = "ol-test0-dede"; TargetCommit = "be90939c9de779fd82751107be1eb50673a8d746"; Path = "e2e/76a5350b-de9e-4e6b-acfd-af5a4cdc3ec6.md" }
)

foreach ($info in $deDeRows) {
    $r = $info.Row
    $wsDeDe.Cells.Item($r, 3).Value = $handedBack               # C: Status
    $xlf = $wsDeDe.Cells.Item($r, 7).Value()                    # G: Latest Handoff File
    $wsDeDe.Cells.Item($r, 9).Value = $info.File                # I: Latest Target File
    $wsDeDe.Cells.Item($r, 10).Value = $xlf                     # J: Latest Handback File
    $wsDeDe.Cells.Item($r, 11).Value = "2016-08-17 06:26:34"    # K: Latest Handback DateTime

    $target = "https://github.com/OpenLocalizationTestOrg/" + $info.TargetRepo + "/blob/" + $info.TargetCommit + "/" + $info.Path
    $wsDeDe.Hyperlinks.Add($wsDeDe.Cells.Item($r, 9), $target, "", "", $info.File)
}
